$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 166666770
$ws.Range("I9").Value = 166666770
$ws.Range("K9").Value = 166666770
$ws.Range("M9").Value = -166666601
$ws.Range("H32").Value = 33830.902
$ws.Range("I32").Value = 6927.375
$ws.Range("J32").Value = 51049.16
$ws.Range("K32").Value = 6927.375
$ws.Range("L32").Value = 51049.16
$ws.Range("M32").Value = -6601.375
$ws.Range("N32").Value = -51701.16
$ws.Range("H86").Value = 5859.5386
$ws.Range("I86").Value = 3862.875
$ws.Range("K86").Value = 3862.875
$ws.Range("M86").Value = -2739.875
$ws.Range("H89").Value = 5859.5386
$ws.Range("I89").Value = 3862.875
$ws.Range("K89").Value = 19314.375
$ws.Range("M89").Value = -13698.375
$ws.Range("H100").Value = 1190.5
$ws.Range("I100").Value = 632.1429000000001
$ws.Range("K100").Value = 632.1429000000001
$ws.Range("M100").Value = -91.14290000000005
$ws.Range("H106").Value = 4931.375
$ws.Range("I106").Value = 4291.1816
$ws.Range("J106").Value = 6339.8
$ws.Range("K106").Value = 4291.1816
$ws.Range("L106").Value = 6339.8
$ws.Range("M106").Value = -3660.1816
$ws.Range("N106").Value = -7601.8
$ws.Range("H133").Value = 95435
$ws.Range("J133").Value = 95435
$ws.Range("L133").Value = 95435
$ws.Range("N133").Value = -105555
$ws.Range("H136").Value = 68169.8
$ws.Range("J136").Value = 81956.86
$ws.Range("L136").Value = 81956.86
$ws.Range("N136").Value = -92156.86
$ws.Range("H137").Value = 350332.97
$ws.Range("I137").Value = 1674.3103
$ws.Range("J137").Value = 1192924.8
$ws.Range("K137").Value = 5022.9309
$ws.Range("L137").Value = 3578774.4
$ws.Range("M137").Value = -2472.9309
$ws.Range("N137").Value = -3583874.4
$ws.Range("H138").Value = 2904.1833
$ws.Range("I138").Value = 2123
$ws.Range("J138").Value = 3356.4473
$ws.Range("K138").Value = 6369
$ws.Range("L138").Value = 10069.3419
$ws.Range("M138").Value = -1229
$ws.Range("N138").Value = -20349.3419

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 77393
$ws.Range("J7").Value = 77393
$ws.Range("L7").Value = 77393
$ws.Range("N7").Value = -77621
$ws.Range("H110").Value = 1441.2727
$ws.Range("I110").Value = 1283.8334
$ws.Range("J110").Value = 2149.75
$ws.Range("K110").Value = 1283.8334
$ws.Range("L110").Value = 2149.75
$ws.Range("M110").Value = 761.1666
$ws.Range("N110").Value = -6239.75
$ws.Range("H117").Value = 73325.2
$ws.Range("J117").Value = 73325.2
$ws.Range("L117").Value = 73325.2
$ws.Range("N117").Value = -82503.2
$ws.Range("H118").Value = 75822.57000000001
$ws.Range("J118").Value = 75822.57000000001
$ws.Range("L118").Value = 75822.57000000001
$ws.Range("N118").Value = -79136.57000000001
$ws.Range("H121").Value = 67306.25
$ws.Range("J121").Value = 67306.25
$ws.Range("L121").Value = 67306.25
$ws.Range("N121").Value = -70800.25
$ws.Range("H132").Value = 3952.8572
$ws.Range("I132").Value = 3884
$ws.Range("J132").Value = 4125
$ws.Range("K132").Value = 11652
$ws.Range("L132").Value = 12375
$ws.Range("M132").Value = -9122
$ws.Range("N132").Value = -17435

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 87400
$ws.Range("J52").Value = 87400
$ws.Range("L52").Value = 87400
$ws.Range("N52").Value = -87926
$ws.Range("H55").Value = 65563
$ws.Range("H86").Value = 2577.375
$ws.Range("I86").Value = 1528
$ws.Range("J86").Value = 3626.75
$ws.Range("K86").Value = 1528
$ws.Range("L86").Value = 3626.75
$ws.Range("M86").Value = -405
$ws.Range("N86").Value = -5872.75
$ws.Range("H89").Value = 2577.375
$ws.Range("I89").Value = 1528
$ws.Range("J89").Value = 3626.75
$ws.Range("K89").Value = 7640
$ws.Range("L89").Value = 18133.75
$ws.Range("M89").Value = -2024
$ws.Range("N89").Value = -29365.75
$ws.Range("H105").Value = 2162.7778
$ws.Range("I105").Value = 2183.125
$ws.Range("K105").Value = 2183.125
$ws.Range("M105").Value = -436.125
$ws.Range("H117").Value = 88086.42999999999
$ws.Range("J117").Value = 88086.42999999999
$ws.Range("L117").Value = 88086.42999999999
$ws.Range("N117").Value = -97264.42999999999
$ws.Range("H121").Value = 87400
$ws.Range("J121").Value = 87400
$ws.Range("L121").Value = 87400
$ws.Range("N121").Value = -90894
$ws.Range("H127").Value = 69880
$ws.Range("J127").Value = 69880
$ws.Range("L127").Value = 69880
$ws.Range("N127").Value = -79800
$ws.Range("H134").Value = 2336.4
$ws.Range("I134").Value = 2204.5881
$ws.Range("K134").Value = 6613.7643
$ws.Range("M134").Value = -4078.7643
$ws.Range("H138").Value = 76478
$ws.Range("J138").Value = 76478
$ws.Range("L138").Value = 76478
$ws.Range("N138").Value = -86758
$ws.Range("H139").Value = 54000
$ws.Range("J139").Value = 54000
$ws.Range("L139").Value = 54000
$ws.Range("N139").Value = -64280
$ws.Range("H140").Value = 58545.75
$ws.Range("J140").Value = 58545.75
$ws.Range("L140").Value = 58545.75
$ws.Range("N140").Value = -68905.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4882.633
$ws.Range("I31").Value = 2853.4211
$ws.Range("K31").Value = 2853.4211
$ws.Range("M31").Value = -2558.4211
$ws.Range("H34").Value = 4882.633
$ws.Range("I34").Value = 2853.4211
$ws.Range("K34").Value = 2853.4211
$ws.Range("M34").Value = -2651.4211
$ws.Range("H105").Value = 87799.69500000001
$ws.Range("I105").Value = 187101.33
$ws.Range("J105").Value = 2684
$ws.Range("K105").Value = 187101.33
$ws.Range("L105").Value = 2684
$ws.Range("M105").Value = -185354.33
$ws.Range("N105").Value = -6178
$ws.Range("H108").Value = 56118.453
$ws.Range("J108").Value = 56118.453
$ws.Range("L108").Value = 56118.453
$ws.Range("N108").Value = -63798.453
$ws.Range("H118").Value = 86098.86
$ws.Range("J118").Value = 86098.86
$ws.Range("L118").Value = 86098.86
$ws.Range("N118").Value = -89412.86
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060
$ws.Range("H141").Value = 152594.03
$ws.Range("J141").Value = 156881.84
$ws.Range("L141").Value = 156881.84
$ws.Range("N141").Value = -167241.84

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 77298.664
$ws.Range("J37").Value = 77298.664
$ws.Range("L37").Value = 231895.992
$ws.Range("N37").Value = -232119.992
$ws.Range("H109").Value = 50002236
$ws.Range("H132").Value = 2528.0967
$ws.Range("I132").Value = 1819.8572
$ws.Range("J132").Value = 2734.6667
$ws.Range("K132").Value = 16378.7148
$ws.Range("L132").Value = 24612.0003
$ws.Range("M132").Value = -13848.7148
$ws.Range("N132").Value = -29672.0003
$ws.Range("H140").Value = 1763.2858
$ws.Range("I140").Value = 1188.6666
$ws.Range("K140").Value = 3565.9998
$ws.Range("M140").Value = 1614.0002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 98323.164
$ws.Range("J114").Value = 98323.164
$ws.Range("L114").Value = 98323.164
$ws.Range("N114").Value = -107001.164
$ws.Range("H116").Value = 50056.47
$ws.Range("J116").Value = 50997.5
$ws.Range("L116").Value = 50997.5
$ws.Range("N116").Value = -60175.5
$ws.Range("H119").Value = 69577.375
$ws.Range("J119").Value = 69577.375
$ws.Range("L119").Value = 69577.375
$ws.Range("N119").Value = -79253.375
$ws.Range("H132").Value = 2505.8
$ws.Range("I132").Value = 2101.125
$ws.Range("J132").Value = 4124.5
$ws.Range("K132").Value = 6303.375
$ws.Range("L132").Value = 12373.5
$ws.Range("M132").Value = -3773.375
$ws.Range("N132").Value = -17433.5
$ws.Range("H135").Value = 68993
$ws.Range("J135").Value = 68993
$ws.Range("L135").Value = 68993
$ws.Range("N135").Value = -79133
$ws.Range("H140").Value = 44298
$ws.Range("J140").Value = 51994.453
$ws.Range("L140").Value = 51994.453
$ws.Range("N140").Value = -62354.453

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 40852
$ws.Range("J121").Value = 41195.715
$ws.Range("L121").Value = 41195.715
$ws.Range("N121").Value = -44689.715

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 7176.316
$ws.Range("I107").Value = 11526.6
$ws.Range("J107").Value = 2342.6667
$ws.Range("K107").Value = 34579.8
$ws.Range("L107").Value = 7028.000100000001
$ws.Range("M107").Value = -32659.8
$ws.Range("N107").Value = -10868.0001
$ws.Range("H122").Value = 4391.129
$ws.Range("I122").Value = 3556.1875
$ws.Range("K122").Value = 10668.5625
$ws.Range("M122").Value = -8218.5625
$ws.Range("H126").Value = 81239.5
$ws.Range("I126").Value = 102614
$ws.Range("K126").Value = 307842
$ws.Range("M126").Value = -305372
